$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append 13 rows (22-34) to the community-smells dataset: 3 new samples for
# Hironsan/anago (ids 21-23) followed by 10 new samples for tensorflow/ranking
# (ids 24-33), continuing straight on from the existing 20 rows.
#
# Row 21 (the last existing row) is reused as a template for both look and
# type: a Paste Values from A21:O21 clones each column's existing cell type --
# numeric id in A, literal text "0"/"1" flags in F:O -- and a separate Paste
# Formats from A21 keeps the new id cells boxed/bold like the rest of the
# table. Typing a fresh "0"/"1" (or a dd/mm/yyyy-shaped date) straight into a
# General cell would make Excel re-sniff it back into a number/date, so any
# cell whose value actually differs from the template is instead overwritten
# via a values-only copy from a cell that already holds that literal text,
# which preserves the text type.
# ---------------------------------------------------------------------------

# --- Row 22 (id 21) ---
$ws.Range("A21:O21").Copy()
$ws.Range("A22:O22").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "https://github.com/Hironsan/anago"
$ws.Range("C22").Value = "anago"
$ws.Range("D22").Value = "Hironsan"
$ws.Range("I21").Copy()
$ws.Range("H22").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("F21").Copy()
$ws.Range("N22").PasteSpecial(-4163)  # xlPasteValues

# --- Row 23 (id 22) ---
$ws.Range("A21:O21").Copy()
$ws.Range("A23:O23").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("A21").Copy()
$ws.Range("A23").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "https://github.com/Hironsan/anago"
$ws.Range("C23").Value = "anago"
$ws.Range("D23").Value = "Hironsan"
$ws.Range("I21").Copy()
$ws.Range("H23").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("F21").Copy()
$ws.Range("N23").PasteSpecial(-4163)  # xlPasteValues

# --- Row 24 (id 23) ---
$ws.Range("A21:O21").Copy()
$ws.Range("A24:O24").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("A21").Copy()
$ws.Range("A24").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "https://github.com/Hironsan/anago"
$ws.Range("C24").Value = "anago"
$ws.Range("D24").Value = "Hironsan"
$ws.Range("I21").Copy()
$ws.Range("H24").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("F21").Copy()
$ws.Range("N24").PasteSpecial(-4163)  # xlPasteValues

# --- Row 25 (id 24) ---
$ws.Range("A21:O21").Copy()
$ws.Range("A25:O25").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("A21").Copy()
$ws.Range("A25").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A25").Value = 24

# --- Row 26 (id 25) ---
$ws.Range("A21:O21").Copy()
$ws.Range("A26:O26").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("A21").Copy()
$ws.Range("A26").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A26").Value = 25

# --- Row 27 (id 26) ---
$ws.Range("A21:O21").Copy()
$ws.Range("A27:O27").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("A21").Copy()
$ws.Range("A27").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A27").Value = 26

# --- Row 28 (id 27) ---
$ws.Range("A21:O21").Copy()
$ws.Range("A28:O28").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("A21").Copy()
$ws.Range("A28").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A28").Value = 27

# --- Row 29 (id 28) ---
$ws.Range("A21:O21").Copy()
$ws.Range("A29:O29").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("A21").Copy()
$ws.Range("A29").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A29").Value = 28
$ws.Range("I21").Copy()
$ws.Range("J29").PasteSpecial(-4163)  # xlPasteValues

# --- Row 30 (id 29) ---
$ws.Range("A21:O21").Copy()
$ws.Range("A30:O30").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("A21").Copy()
$ws.Range("A30").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A30").Value = 29

# --- Row 31 (id 30) ---
$ws.Range("A21:O21").Copy()
$ws.Range("A31:O31").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("A21").Copy()
$ws.Range("A31").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A31").Value = 30
$ws.Range("I21").Copy()
$ws.Range("J31").PasteSpecial(-4163)  # xlPasteValues

# --- Row 32 (id 31) ---
$ws.Range("A21:O21").Copy()
$ws.Range("A32:O32").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("A21").Copy()
$ws.Range("A32").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A32").Value = 31

# --- Row 33 (id 32) ---
$ws.Range("A21:O21").Copy()
$ws.Range("A33:O33").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("A21").Copy()
$ws.Range("A33").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A33").Value = 32
$ws.Range("I21").Copy()
$ws.Range("J33").PasteSpecial(-4163)  # xlPasteValues

# --- Row 34 (id 33) ---
$ws.Range("A21:O21").Copy()
$ws.Range("A34:O34").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("A21").Copy()
$ws.Range("A34").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A34").Value = 33

# ---------------------------------------------------------------------------
# Dates that differ from the row-21 template ("12/03/2018") are patched in
# last, via an off-sheet helper cell that mints the literal string as text
# through a TEXT() formula (so Excel never gets a chance to parse it back into
# a date), Paste-Valued into place, then torn down again.
# ---------------------------------------------------------------------------
$helper = $ws.Range("ZZ1")
$helper.Formula = "=TEXT(DATE(2017,6,26),""mm/dd/yyyy"")"
$helper.Copy()
$ws.Range("E22").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E23").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E24").PasteSpecial(-4163)  # xlPasteValues

$helper.Clear()
$ws.Columns("ZZ").Delete()

$ws.Application.CutCopyMode = $false

